# Updated cryptos list - applies per-cell value changes from the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.644.19'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '3.675.25'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.76'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('E6').Value = '  +8.68%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '661.17'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +1.44%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = '3.673.22'
$ws.Range('E11').Value = '  +1.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.89'
$ws.Range('E12').Value = '  +4.47%  '
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('E14').Value = '  +4.62%  '
$ws.Range('D15').Value = '4.360.30'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000268'
$ws.Range('E16').Value = '  +2.60%  '
$ws.Range('D17').Value = '96.464.66'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.91'
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('D19').Value = '3.680.65'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.89'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.71'
$ws.Range('E21').Value = '  +2.89%  '
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '525.65'
$ws.Range('E23').Value = '  +2.32%  '
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('E25').Value = '  +2.66%  '
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.24'
$ws.Range('E27').Value = '  +3.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.16'
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').Value = '3.873.49'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.169'
$ws.Range('E30').Value = '  +10.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.54'
$ws.Range('E31').Value = '  +6.44%  '
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('E34').Value = '  +15.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.187'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('B36').Value = 'Bittensor'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '664.47'
$ws.Range('E36').Value = '  +7.14%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '32.59'
$ws.Range('E37').Value = '  +2.41%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.592'
$ws.Range('E39').Value = '  +3.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.89'
$ws.Range('E40').Value = '  +1.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '43.10'
$ws.Range('E41').Value = '  +29.06%  '
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.99'
$ws.Range('E43').Value = '  +1.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.966'
$ws.Range('E44').Value = '  +3.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.55'
$ws.Range('E45').Value = '  +9.71%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  +6.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.448'
$ws.Range('E48').Value = '  +16.34%  '
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('E51').Value = '  +1.22%  '
